$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CRM 3586 - Remove the "Service Tax No." column (header in row 1 /
# {vendor:service_tax_no} merge-field in row 2) from the SF List template.
# This is column R; deleting it shifts every following column (S..AF) one
# place to the left, which also drops the two now-unused shared-string
# entries ("Service Tax No." and "{vendor:service_tax_no}").
$ws.Columns("R").Delete()

# Reflect the selection left behind after the column removal: the user's
# cursor ends up on the newly-trailing (now empty) column AF.
$ws.Range("AF1:AF1048576").Select()
